$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = 0.239058013728998
$ws.Range("E3").Value = 0.4242326590369501
$ws.Range("D4").Value = 0.08141511437990984
$ws.Range("D6").Value = 0.01449024281269066
$ws.Range("D7").Value = 0.02475059564391487
$ws.Range("D9").Value = 0.04134373707948353
$ws.Range("E9").Value = 0.04784006224598872
$ws.Range("D11").Value = 0.285210188799379
$ws.Range("E11").Value = 0.3235559504570613
$ws.Range("F11").Value = 0.3407827236925464
$ws.Range("D12").Value = 0.01287395819840383
$ws.Range("D13").Value = 3.943365448644887
$ws.Range("E13").Value = 5.402881633126435
$ws.Range("F13").Value = 9.125743175701356
$ws.Range("D15").Value = 0.2288425564963268
$ws.Range("D16").Value = 0.1699639261961932
$ws.Range("D19").Value = 0.02408145577782827
$ws.Range("D23").Value = 0.4539387891005192
$ws.Range("D24").Value = 0.02203609467763762
$ws.Range("D25").Value = 11.0041087037838
$ws.Range("E25").Value = 14.99176931333146
$ws.Range("F25").Value = 17.699528219766
$ws.Range("D28").Value = 0.1604933318719226
$ws.Range("D30").Value = 0.01461259042098011
$ws.Range("D31").Value = 0.02292538174273028
$ws.Range("E31").Value = 0.3367918963752413
$ws.Range("D33").Value = 0.04583697821204784
$ws.Range("D35").Value = 0.3150216178095544
$ws.Range("E35").Value = 0.3432368386291454
$ws.Range("D37").Value = 14.1019951135758
$ws.Range("E37").Value = 21.10008000308873
$ws.Range("F37").Value = 25.73458763814209
$ws.Range("D40").Value = 0.08009540424700452
$ws.Range("E40").Value = 0.08456327765055979
$ws.Range("F40").Value = 0.2295091162223122
$ws.Range("D43").Value = 0.02390075507836695
$ws.Range("E43").Value = 0.02475138362629163
$ws.Range("D44").Value = 0.2392629983996759
$ws.Range("D45").Value = 0.122677231353464
$ws.Range("D48").Value = 0.01463024846736994
$ws.Range("E48").Value = 0.0148620061315173
$ws.Range("D49").Value = 14.32196783451306
$ws.Range("E49").Value = 17.16468029822554
$ws.Range("F49").Value = 25.79560824674561
